$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 14:10"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1746217
$ws.Range("C4").Value = 414
$ws.Range("D4").Value = 490151
$ws.Range("E4").Value = 1153952
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 102114

# Row 13: India
$ws.Range("A13").Value = "India"
$ws.Range("B13").Value = 159054
$ws.Range("C13").Value = 968
$ws.Range("D13").Value = 67929
$ws.Range("E13").Value = 86584
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 7
$ws.Range("H13").Value = 4541

# Row 28: Suecia
$ws.Range("A28").Value = "Suecia"
$ws.Range("B28").Value = 35727
$ws.Range("C28").Value = 639
$ws.Range("D28").Value = 4971
$ws.Range("E28").Value = 26490
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 46
$ws.Range("H28").Value = 4266

# Row 50: Dinamarca
$ws.Range("A50").Value = "Dinamarca"
$ws.Range("B50").Value = 11512
$ws.Range("C50").Value = 32
$ws.Range("D50").Value = 10180
$ws.Range("E50").Value = 764
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 568

# Row 76: Uzbekistan
$ws.Range("A76").Value = "Uzbekistan"
$ws.Range("B76").Value = 3437
$ws.Range("C76").Value = 68
$ws.Range("D76").Value = 2685
$ws.Range("E76").Value = 738
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 14

# Row 77: Tayikistan
$ws.Range("A77").Value = "Tayikistan"
$ws.Range("B77").Value = 3424
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 1575
$ws.Range("E77").Value = 1802
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 47

# Row 107: Libano
$ws.Range("A107").Value = "Libano"
$ws.Range("B107").Value = 1168
$ws.Range("C107").Value = 7
$ws.Range("D107").Value = 699
$ws.Range("E107").Value = 443
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 26

# Row 190: San Vicente y las Granadinas
$ws.Range("A190").Value = "San Vicente y las Granadinas"
$ws.Range("B190").Value = 25
$ws.Range("C190").Value = 7
$ws.Range("D190").Value = 14
$ws.Range("E190").Value = 11
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

# Row 191: Gambia
$ws.Range("A191").Value = "Gambia"
$ws.Range("B191").Value = 25
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 19
$ws.Range("E191").Value = 5
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 1

# Row 192: Antigua y Barbuda
$ws.Range("A192").Value = "Antigua y Barbuda"
$ws.Range("B192").Value = 25
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 19
$ws.Range("E192").Value = 3
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 3

# Row 193: Timor Oriental
$ws.Range("A193").Value = "Timor Oriental"
$ws.Range("B193").Value = 24
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 24
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

# Row 194: Granada
$ws.Range("A194").Value = "Granada"
$ws.Range("B194").Value = 23
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 18
$ws.Range("E194").Value = 5
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0

# Row 195: Namibia
$ws.Range("A195").Value = "Namibia"
$ws.Range("B195").Value = 22
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 14
$ws.Range("E195").Value = 8
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0

# Row 196: Laos
$ws.Range("A196").Value = "Laos"
$ws.Range("B196").Value = 19
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 16
$ws.Range("E196").Value = 3
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

# Row 197: Fiyi
$ws.Range("A197").Value = "Fiyi"
$ws.Range("B197").Value = 18
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 15
$ws.Range("E197").Value = 3
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 0

# Row 198: Curazao
$ws.Range("A198").Value = "Curazao"
$ws.Range("B198").Value = 18
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 14
$ws.Range("E198").Value = 3
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

# Row 199: Nueva Caledonia
$ws.Range("A199").Value = "Nueva Caledonia"
$ws.Range("B199").Value = 18
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 18
$ws.Range("E199").Value = 0
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

# Row 200: Santa Lucia
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 18
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0

# Row 201: Belice
$ws.Range("A201").Value = "Belice"
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 16
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 2

# Row 210: Seychelles
$ws.Range("A210").Value = "Seychelles"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 11
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211: Montserrat
$ws.Range("A211").Value = "Montserrat"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 10
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 1

# Row 213: Islas Virgenes Britanicas
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# Row 214: Papua Nueva Guinea
$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# Row 215: Bonaire, San Eustaquio y Saba
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B215").Value = 6
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 6
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

# Row 216: San Bartolome
$ws.Range("A216").Value = "San Bartolome"
$ws.Range("B216").Value = 6
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 6
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0

